$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$style = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '92.297.27'
$ws.Range('D2').Style = $style
$ws.Range('E2').Value = '  +0.80%  '

# Row 3
$style = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.099.29'
$ws.Range('D3').Style = $style
$ws.Range('E3').Value = '  -1.53%  '

# Row 4
$ws.Range('E4').Value = '  -0.08%  '

# Row 5
$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.03'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  -2.78%  '

# Row 6
$style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '613.22'
$ws.Range('D6').Style = $style
$ws.Range('E6').Value = '  -0.87%  '

# Row 7
$style = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.09'
$ws.Range('D7').Style = $style
$ws.Range('E7').Value = '  -2.47%  '

# Row 8
$style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.387'
$ws.Range('D8').Style = $style
$ws.Range('E8').Value = '  -0.16%  '

# Row 9
$ws.Range('E9').Value = '  -0.09%  '

# Row 10
$style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.094.91'
$ws.Range('D10').Style = $style
$ws.Range('E10').Value = '  -1.68%  '

# Row 11
$style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.781'
$ws.Range('D11').Style = $style
$ws.Range('E11').Value = '  +5.42%  '

# Row 12
$ws.Range('E12').Value = '  -3.28%  '

# Row 13
$ws.Range('E13').Value = '  -4.49%  '

# Row 14
$style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '92.024.93'
$ws.Range('D14').Style = $style
$ws.Range('E14').Value = '  +0.83%  '

# Row 15
$style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '33.79'
$ws.Range('D15').Style = $style
$ws.Range('E15').Value = '  -3.30%  '

# Row 16
$style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.40'
$ws.Range('D16').Style = $style
$ws.Range('E16').Value = '  -3.34%  '

# Row 18
$ws.Range('E18').Value = '  -2.74%  '

# Row 19
$ws.Range('E19').Value = '  +1.72%  '

# Row 20
$style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.39'
$ws.Range('D20').Style = $style
$ws.Range('E20').Value = '  -3.57%  '

# Row 21
$ws.Range('E21').Value = '  -2.23%  '

# Row 22
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '436.26'
$ws.Range('D22').Style = $style
$ws.Range('E22').Value = '  -4.35%  '

# Row 23
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.09'
$ws.Range('D23').Style = $style
$ws.Range('E23').Value = '  -0.90%  '

# Row 24
$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000198'
$ws.Range('D24').Style = $style
$ws.Range('E24').Value = '  -2.00%  '

# Row 25
$style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.56'
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  -5.38%  '

# Row 26
$style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '85.27'
$ws.Range('D26').Style = $style
$ws.Range('E26').Value = '  -3.70%  '

# Row 27
$style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.35'
$ws.Range('D27').Style = $style
$ws.Range('E27').Value = '  -3.41%  '

# Row 28
$style = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.259.31'
$ws.Range('D28').Style = $style

# Row 30
$ws.Range('E30').Value = '  +6.68%  '

# Row 31
$style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.235'
$ws.Range('D31').Style = $style
$ws.Range('E31').Value = '  +2.83%  '

# Row 32
$style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.123'
$ws.Range('D32').Style = $style
$ws.Range('E32').Value = '  -18.89%  '

# Row 33
$ws.Range('E33').Value = '  -35.72%  '

# Row 34
$style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.13'
$ws.Range('D34').Style = $style
$ws.Range('E34').Value = '  -2.52%  '

# Row 35
$style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '8.02'
$ws.Range('D35').Style = $style
$ws.Range('E35').Value = '  +7.74%  '

# Row 36
$ws.Range('E36').Value = '  -11.63%  '

# Row 37
$style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '25.57'
$ws.Range('D37').Style = $style
$ws.Range('E37').Value = '  -2.77%  '

# Row 38
$style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.89'
$ws.Range('D38').Style = $style
$ws.Range('E38').Value = '  -0.81%  '

# Row 39
$style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.89'
$ws.Range('D39').Style = $style
$ws.Range('E39').Value = '  -6.16%  '

# Row 40
$style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.87'
$ws.Range('D40').Style = $style
$ws.Range('E40').Value = '  +7.79%  '

# Row 41
$style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.436'
$ws.Range('D41').Style = $style
$ws.Range('E41').Value = '  -1.07%  '

# Row 42
$ws.Range('E42').Value = '  -3.70%  '

# Row 43
$style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '465.36'
$ws.Range('D43').Style = $style
$ws.Range('E43').Value = '  -5.32%  '

# Row 44
$style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.26'
$ws.Range('D44').Style = $style
$ws.Range('E44').Value = '  -3.81%  '

# Row 45
$ws.Range('E45').Value = '  +0.13%  '

# Row 46
$style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '159.03'
$ws.Range('D46').Style = $style
$ws.Range('E46').Value = '  +2.10%  '

# Row 47
$style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.680'
$ws.Range('D47').Style = $style
$ws.Range('E47').Value = '  -3.73%  '

# Row 48
$style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.83'
$ws.Range('D48').Style = $style
$ws.Range('E48').Value = '  -4.97%  '

# Row 49
$style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '43.74'
$ws.Range('D49').Style = $style
$ws.Range('E49').Value = '  -0.73%  '

# Row 50
$ws.Range('E50').Value = '  -3.20%  '

# Row 51
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$style = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0323'
$ws.Range('D51').Style = $style
$ws.Range('E51').Value = '  +0.99%  '
